# Auto-generated edit script: updates LeveProfits-style price/profit
# columns (H-N) across multiple sheets per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 426.33334
$ws.Cells.Item(18, 9).Value = 414.5
$ws.Cells.Item(18, 11).Value = 414.5
$ws.Cells.Item(18, 13).Value = -130.5
$ws.Cells.Item(40, 8).Value = 9832.333000000001
$ws.Cells.Item(40, 9).Value = 4664.6665
$ws.Cells.Item(40, 10).Value = 15000
$ws.Cells.Item(40, 11).Value = 4664.6665
$ws.Cells.Item(40, 12).Value = 15000
$ws.Cells.Item(40, 13).Value = -4489.6665
$ws.Cells.Item(40, 14).Value = -15350
$ws.Cells.Item(51, 8).Value = 6493.75
$ws.Cells.Item(51, 9).Value = 5000
$ws.Cells.Item(51, 10).Value = 6991.6665
$ws.Cells.Item(51, 11).Value = 5000
$ws.Cells.Item(51, 12).Value = 6991.6665
$ws.Cells.Item(51, 13).Value = -4516
$ws.Cells.Item(51, 14).Value = -7959.6665
$ws.Cells.Item(64, 8).Value = 8011
$ws.Cells.Item(64, 9).Value = 6522
$ws.Cells.Item(64, 10).Value = 9500
$ws.Cells.Item(64, 11).Value = 6522
$ws.Cells.Item(64, 12).Value = 9500
$ws.Cells.Item(64, 13).Value = -6274
$ws.Cells.Item(64, 14).Value = -9996
$ws.Cells.Item(67, 8).Value = 8011
$ws.Cells.Item(67, 9).Value = 6522
$ws.Cells.Item(67, 10).Value = 9500
$ws.Cells.Item(67, 11).Value = 6522
$ws.Cells.Item(67, 12).Value = 9500
$ws.Cells.Item(67, 13).Value = -5664
$ws.Cells.Item(67, 14).Value = -11216
$ws.Cells.Item(69, 8).Value = 7347.143
$ws.Cells.Item(69, 10).Value = 7905
$ws.Cells.Item(69, 12).Value = 23715
$ws.Cells.Item(69, 14).Value = -25463
$ws.Cells.Item(72, 8).Value = 7347.143
$ws.Cells.Item(72, 10).Value = 7905
$ws.Cells.Item(72, 12).Value = 71145
$ws.Cells.Item(72, 14).Value = -79881
$ws.Cells.Item(116, 8).Value = 3975.7144
$ws.Cells.Item(116, 9).Value = 3388.3333
$ws.Cells.Item(116, 11).Value = 3388.3333
$ws.Cells.Item(116, 13).Value = 53.66670000000022
$ws.Cells.Item(125, 8).Value = 2488.9285
$ws.Cells.Item(125, 9).Value = 1088.2858
$ws.Cells.Item(125, 11).Value = 9794.572200000001
$ws.Cells.Item(125, 13).Value = -7334.572200000001
$ws.Cells.Item(137, 8).Value = 3090.2173
$ws.Cells.Item(137, 9).Value = 1876.909
$ws.Cells.Item(137, 10).Value = 4202.4165
$ws.Cells.Item(137, 11).Value = 5630.727000000001
$ws.Cells.Item(137, 12).Value = 12607.2495
$ws.Cells.Item(137, 13).Value = -3080.727000000001
$ws.Cells.Item(137, 14).Value = -17707.2495
$ws.Cells.Item(138, 8).Value = 3034.48
$ws.Cells.Item(138, 10).Value = 3368.7812
$ws.Cells.Item(138, 12).Value = 10106.3436
$ws.Cells.Item(138, 14).Value = -20386.3436

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(11, 8).Value = 176.25
$ws.Cells.Item(11, 9).Value = 127.5
$ws.Cells.Item(11, 11).Value = 127.5
$ws.Cells.Item(11, 13).Value = 16.5
$ws.Cells.Item(32, 8).Value = 15588.625
$ws.Cells.Item(32, 10).Value = 45832.777
$ws.Cells.Item(32, 12).Value = 45832.777
$ws.Cells.Item(32, 14).Value = -46406.777
$ws.Cells.Item(61, 8).Value = 3921.7
$ws.Cells.Item(61, 9).Value = 2182.6
$ws.Cells.Item(61, 10).Value = 7399.9
$ws.Cells.Item(61, 11).Value = 2182.6
$ws.Cells.Item(61, 12).Value = 7399.9
$ws.Cells.Item(61, 13).Value = -1970.6
$ws.Cells.Item(61, 14).Value = -7823.9
$ws.Cells.Item(136, 8).Value = 3921.7
$ws.Cells.Item(136, 9).Value = 2182.6
$ws.Cells.Item(136, 10).Value = 7399.9
$ws.Cells.Item(136, 11).Value = 6547.799999999999
$ws.Cells.Item(136, 12).Value = 22199.7
$ws.Cells.Item(136, 13).Value = -3997.799999999999
$ws.Cells.Item(136, 14).Value = -27299.7

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3858.8484
$ws.Cells.Item(134, 9).Value = 2855.5715
$ws.Cells.Item(134, 11).Value = 8566.7145
$ws.Cells.Item(134, 13).Value = -6031.7145

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 753.4211
$ws.Cells.Item(22, 9).Value = 616.61536
$ws.Cells.Item(22, 11).Value = 616.61536
$ws.Cells.Item(22, 13).Value = -266.61536
$ws.Cells.Item(31, 8).Value = 3756.946
$ws.Cells.Item(31, 9).Value = 2627.7222
$ws.Cells.Item(31, 11).Value = 2627.7222
$ws.Cells.Item(31, 13).Value = -2332.7222
$ws.Cells.Item(34, 8).Value = 3756.946
$ws.Cells.Item(34, 9).Value = 2627.7222
$ws.Cells.Item(34, 11).Value = 2627.7222
$ws.Cells.Item(34, 13).Value = -2425.7222
$ws.Cells.Item(58, 8).Value = 3125
$ws.Cells.Item(58, 9).Value = 1833.3334
$ws.Cells.Item(58, 10).Value = 7000
$ws.Cells.Item(58, 11).Value = 1833.3334
$ws.Cells.Item(58, 12).Value = 7000
$ws.Cells.Item(58, 13).Value = -1630.3334
$ws.Cells.Item(58, 14).Value = -7406
$ws.Cells.Item(122, 8).Value = 271694
$ws.Cells.Item(122, 9).Value = 465189.9
$ws.Cells.Item(122, 10).Value = 5637.125
$ws.Cells.Item(122, 11).Value = 1395569.7
$ws.Cells.Item(122, 12).Value = 16911.375
$ws.Cells.Item(122, 13).Value = -1393119.7
$ws.Cells.Item(122, 14).Value = -21811.375
$ws.Cells.Item(136, 8).Value = 3125
$ws.Cells.Item(136, 9).Value = 1833.3334
$ws.Cells.Item(136, 10).Value = 7000
$ws.Cells.Item(136, 11).Value = 5500.0002
$ws.Cells.Item(136, 12).Value = 21000
$ws.Cells.Item(136, 13).Value = -2950.0002
$ws.Cells.Item(136, 14).Value = -26100

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 269
$ws.Cells.Item(2, 10).Value = 329.41666
$ws.Cells.Item(2, 12).Value = 1976.49996
$ws.Cells.Item(2, 14).Value = -2202.49996
$ws.Cells.Item(9, 8).Value = 909182.4399999999
$ws.Cells.Item(9, 9).Value = 333434.1
$ws.Cells.Item(9, 10).Value = 6666666
$ws.Cells.Item(9, 11).Value = 1000302.3
$ws.Cells.Item(9, 12).Value = 19999998
$ws.Cells.Item(9, 13).Value = -1000078.3
$ws.Cells.Item(9, 14).Value = -20000446
$ws.Cells.Item(34, 8).Value = 2323.182
$ws.Cells.Item(34, 10).Value = 4987.5
$ws.Cells.Item(34, 12).Value = 14962.5
$ws.Cells.Item(34, 14).Value = -15130.5
$ws.Cells.Item(39, 8).Value = 4408.409
$ws.Cells.Item(39, 10).Value = 4408.409
$ws.Cells.Item(39, 12).Value = 13225.227
$ws.Cells.Item(39, 14).Value = -13813.227
$ws.Cells.Item(56, 8).Value = 7666.353
$ws.Cells.Item(56, 9).Value = 7666.353
$ws.Cells.Item(56, 11).Value = 7666.353
$ws.Cells.Item(56, 13).Value = -7136.353
$ws.Cells.Item(80, 8).Value = 8224.625
$ws.Cells.Item(80, 10).Value = 8828.286
$ws.Cells.Item(80, 12).Value = 26484.858
$ws.Cells.Item(80, 14).Value = -28356.858
$ws.Cells.Item(83, 8).Value = 8224.625
$ws.Cells.Item(83, 10).Value = 8828.286
$ws.Cells.Item(83, 12).Value = 79454.57399999999
$ws.Cells.Item(83, 14).Value = -88814.57399999999
$ws.Cells.Item(92, 8).Value = 699.6667
$ws.Cells.Item(92, 9).Value = 900
$ws.Cells.Item(92, 10).Value = 659.6
$ws.Cells.Item(92, 11).Value = 2700
$ws.Cells.Item(92, 12).Value = 1978.8
$ws.Cells.Item(92, 13).Value = -1452
$ws.Cells.Item(92, 14).Value = -4474.8
$ws.Cells.Item(113, 8).Value = 1770.6923
$ws.Cells.Item(113, 9).Value = 1741.6666
$ws.Cells.Item(113, 10).Value = 1779.4
$ws.Cells.Item(113, 11).Value = 5224.9998
$ws.Cells.Item(113, 12).Value = 5338.200000000001
$ws.Cells.Item(113, 13).Value = -3054.9998
$ws.Cells.Item(113, 14).Value = -9678.200000000001
$ws.Cells.Item(122, 8).Value = 1206.2903
$ws.Cells.Item(122, 9).Value = 1201.4286
$ws.Cells.Item(122, 10).Value = 1207.7084
$ws.Cells.Item(122, 11).Value = 10812.8574
$ws.Cells.Item(122, 12).Value = 10869.3756
$ws.Cells.Item(122, 13).Value = -8362.857399999999
$ws.Cells.Item(122, 14).Value = -15769.3756

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 24166.166
$ws.Cells.Item(15, 10).Value = 24000
$ws.Cells.Item(15, 12).Value = 24000
$ws.Cells.Item(15, 14).Value = -24576
$ws.Cells.Item(81, 8).Value = 24166.166
$ws.Cells.Item(81, 10).Value = 24000
$ws.Cells.Item(81, 12).Value = 24000
$ws.Cells.Item(81, 14).Value = -25996
$ws.Cells.Item(84, 8).Value = 24166.166
$ws.Cells.Item(84, 10).Value = 24000
$ws.Cells.Item(84, 12).Value = 72000
$ws.Cells.Item(84, 14).Value = -81984

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2152
$ws.Cells.Item(46, 9).Value = 1574.25
$ws.Cells.Item(46, 10).Value = 3142.4285
$ws.Cells.Item(46, 11).Value = 1574.25
$ws.Cells.Item(46, 12).Value = 3142.4285
$ws.Cells.Item(46, 13).Value = -1386.25
$ws.Cells.Item(46, 14).Value = -3518.4285
$ws.Cells.Item(136, 8).Value = 4895
$ws.Cells.Item(136, 9).Value = 3050.05
$ws.Cells.Item(136, 11).Value = 9150.150000000001
$ws.Cells.Item(136, 13).Value = -6600.150000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3876.0952
$ws.Cells.Item(122, 9).Value = 2880.6
$ws.Cells.Item(122, 11).Value = 8641.799999999999
$ws.Cells.Item(122, 13).Value = -6191.799999999999
$ws.Cells.Item(136, 8).Value = 11908739
$ws.Cells.Item(136, 9).Value = 18521506
$ws.Cells.Item(136, 11).Value = 55564518
$ws.Cells.Item(136, 13).Value = -55561968

Write-Output "Applied 201 cell updates across 8 sheets."
